# Apply updated crypto price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '36.605.54'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -2.16%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.037.64'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  -0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '231.73'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -9.83%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.602'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("E7").Value = '  -0.35%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '55.22'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.74%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.371'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.95%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '56.89'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.06%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0753'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.47%  '

$ws.Range("E12").Value = '  +0.55%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '2.334.72'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.25%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '14.33'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.87%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '20.05'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -8.15%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.760'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '5.17'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.09%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.034.59'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.55%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '36.707.28'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.50%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.84'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +15.28%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '67.57'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.14%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.0₃0797'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.75%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '220.82'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -5.76%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("E25").Value = '  +1.24%  '

$ws.Range("E26").Value = '  -7.26%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '162.86'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.16%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '8.72'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.46%  '

$ws.Range("E29").Value = '  -3.00%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '18.95'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.43%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.35'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +3.40%  '

$ws.Range("E32").Value = '  -1.22%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.36'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.10%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0606'
$c.Style = "Normal"

$ws.Range("E35").Value = '  +3.69%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.27'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.98%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("E38").Value = '  -2.71%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.79'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +7.98%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '3.22'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -5.82%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '4.51'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +41.17%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.95'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -3.56%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.477.11'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.34%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0940'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.84%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '93.56'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +4.73%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0204'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("E47").Value = '  -4.58%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '15.59'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.17%  '

$ws.Range("E49").Value = '  -2.60%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.89'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.17%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.94'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.88%  '
